# Auto update Excel log: append newly-logged sensor readings (2026-02-06)
# to the PIR, Humidity and Temperature sheets.
$wb = $excel.ActiveWorkbook

# Helper: write a value that might look like a date/time/number (e.g. '2026-02-06',
# '10:11:30', '68.8%') while forcing it to be stored as plain text, matching the
# existing log rows (which are all text, not real dates/numbers).
function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Helper: write a value that Excel will never misinterpret (plain words like
# 'Bathroom', 'Active', 'Inactive', 'No Motion', ...).
function Set-PlainCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# PIR sheet - append rows 379-392
$ws = $wb.Worksheets.Item("PIR")
Set-TextCell $ws 379 1 '2026-02-06'
Set-TextCell $ws 379 2 '10:11:30'
Set-TextCell $ws 379 3 '10:00'
Set-PlainCell $ws 379 4 'Bathroom'
Set-PlainCell $ws 379 5 'No Motion'
Set-PlainCell $ws 379 6 'Inactive'

Set-TextCell $ws 380 1 '2026-02-06'
Set-TextCell $ws 380 2 '10:11:33'
Set-TextCell $ws 380 3 '10:00'
Set-PlainCell $ws 380 4 'Bathroom'
Set-PlainCell $ws 380 5 'No Motion'
Set-PlainCell $ws 380 6 'Inactive'

Set-TextCell $ws 381 1 '2026-02-06'
Set-TextCell $ws 381 2 '10:11:35'
Set-TextCell $ws 381 3 '10:00'
Set-PlainCell $ws 381 4 'Bathroom'
Set-PlainCell $ws 381 5 'No Motion'
Set-PlainCell $ws 381 6 'Inactive'

Set-TextCell $ws 382 1 '2026-02-06'
Set-TextCell $ws 382 2 '10:11:39'
Set-TextCell $ws 382 3 '10:00'
Set-PlainCell $ws 382 4 'Bathroom'
Set-PlainCell $ws 382 5 'No Motion'
Set-PlainCell $ws 382 6 'Inactive'

Set-TextCell $ws 383 1 '2026-02-06'
Set-TextCell $ws 383 2 '10:11:44'
Set-TextCell $ws 383 3 '10:00'
Set-PlainCell $ws 383 4 'Bathroom'
Set-PlainCell $ws 383 5 'No Motion'
Set-PlainCell $ws 383 6 'Inactive'

Set-TextCell $ws 384 1 '2026-02-06'
Set-TextCell $ws 384 2 '10:11:50'
Set-TextCell $ws 384 3 '10:00'
Set-PlainCell $ws 384 4 'Bathroom'
Set-PlainCell $ws 384 5 'No Motion'
Set-PlainCell $ws 384 6 'Inactive'

Set-TextCell $ws 385 1 '2026-02-06'
Set-TextCell $ws 385 2 '10:11:50'
Set-TextCell $ws 385 3 '10:00'
Set-PlainCell $ws 385 4 'Bathroom'
Set-PlainCell $ws 385 5 'Motion Detected'
Set-PlainCell $ws 385 6 'Active'

Set-TextCell $ws 386 1 '2026-02-06'
Set-TextCell $ws 386 2 '10:11:56'
Set-TextCell $ws 386 3 '10:00'
Set-PlainCell $ws 386 4 'Bathroom'
Set-PlainCell $ws 386 5 'No Motion'
Set-PlainCell $ws 386 6 'Inactive'

Set-TextCell $ws 387 1 '2026-02-06'
Set-TextCell $ws 387 2 '10:11:59'
Set-TextCell $ws 387 3 '10:00'
Set-PlainCell $ws 387 4 'Bathroom'
Set-PlainCell $ws 387 5 'Motion Detected'
Set-PlainCell $ws 387 6 'Active'

Set-TextCell $ws 388 1 '2026-02-06'
Set-TextCell $ws 388 2 '10:12:05'
Set-TextCell $ws 388 3 '10:00'
Set-PlainCell $ws 388 4 'Bathroom'
Set-PlainCell $ws 388 5 'No Motion'
Set-PlainCell $ws 388 6 'Inactive'

Set-TextCell $ws 389 1 '2026-02-06'
Set-TextCell $ws 389 2 '10:12:10'
Set-TextCell $ws 389 3 '10:00'
Set-PlainCell $ws 389 4 'Bathroom'
Set-PlainCell $ws 389 5 'No Motion'
Set-PlainCell $ws 389 6 'Inactive'

Set-TextCell $ws 390 1 '2026-02-06'
Set-TextCell $ws 390 2 '10:12:15'
Set-TextCell $ws 390 3 '10:00'
Set-PlainCell $ws 390 4 'Bathroom'
Set-PlainCell $ws 390 5 'No Motion'
Set-PlainCell $ws 390 6 'Inactive'

Set-TextCell $ws 391 1 '2026-02-06'
Set-TextCell $ws 391 2 '10:12:20'
Set-TextCell $ws 391 3 '10:00'
Set-PlainCell $ws 391 4 'Bathroom'
Set-PlainCell $ws 391 5 'No Motion'
Set-PlainCell $ws 391 6 'Inactive'

Set-TextCell $ws 392 1 '2026-02-06'
Set-TextCell $ws 392 2 '10:12:25'
Set-TextCell $ws 392 3 '10:00'
Set-PlainCell $ws 392 4 'Bathroom'
Set-PlainCell $ws 392 5 'No Motion'
Set-PlainCell $ws 392 6 'Inactive'

# Humidity sheet - append rows 260-270
$ws = $wb.Worksheets.Item("Humidity")
Set-TextCell $ws 260 1 '2026-02-06'
Set-TextCell $ws 260 2 '10:11:31'
Set-TextCell $ws 260 3 '10:00'
Set-PlainCell $ws 260 4 'Bathroom'
Set-TextCell $ws 260 5 '68.8%'
Set-PlainCell $ws 260 6 'Active'

Set-TextCell $ws 261 1 '2026-02-06'
Set-TextCell $ws 261 2 '10:11:34'
Set-TextCell $ws 261 3 '10:00'
Set-PlainCell $ws 261 4 'Bathroom'
Set-TextCell $ws 261 5 '68.0%'
Set-PlainCell $ws 261 6 'Active'

Set-TextCell $ws 262 1 '2026-02-06'
Set-TextCell $ws 262 2 '10:11:38'
Set-TextCell $ws 262 3 '10:00'
Set-PlainCell $ws 262 4 'Bathroom'
Set-TextCell $ws 262 5 '69.0%'
Set-PlainCell $ws 262 6 'Active'

Set-TextCell $ws 263 1 '2026-02-06'
Set-TextCell $ws 263 2 '10:11:43'
Set-TextCell $ws 263 3 '10:00'
Set-PlainCell $ws 263 4 'Bathroom'
Set-TextCell $ws 263 5 '67.9%'
Set-PlainCell $ws 263 6 'Active'

Set-TextCell $ws 264 1 '2026-02-06'
Set-TextCell $ws 264 2 '10:11:48'
Set-TextCell $ws 264 3 '10:00'
Set-PlainCell $ws 264 4 'Bathroom'
Set-TextCell $ws 264 5 '68.8%'
Set-PlainCell $ws 264 6 'Active'

Set-TextCell $ws 265 1 '2026-02-06'
Set-TextCell $ws 265 2 '10:11:53'
Set-TextCell $ws 265 3 '10:00'
Set-PlainCell $ws 265 4 'Bathroom'
Set-TextCell $ws 265 5 '67.8%'
Set-PlainCell $ws 265 6 'Active'

Set-TextCell $ws 266 1 '2026-02-06'
Set-TextCell $ws 266 2 '10:11:58'
Set-TextCell $ws 266 3 '10:00'
Set-PlainCell $ws 266 4 'Bathroom'
Set-TextCell $ws 266 5 '68.9%'
Set-PlainCell $ws 266 6 'Active'

Set-TextCell $ws 267 1 '2026-02-06'
Set-TextCell $ws 267 2 '10:12:08'
Set-TextCell $ws 267 3 '10:00'
Set-PlainCell $ws 267 4 'Bathroom'
Set-TextCell $ws 267 5 '67.4%'
Set-PlainCell $ws 267 6 'Active'

Set-TextCell $ws 268 1 '2026-02-06'
Set-TextCell $ws 268 2 '10:12:13'
Set-TextCell $ws 268 3 '10:00'
Set-PlainCell $ws 268 4 'Bathroom'
Set-TextCell $ws 268 5 '68.1%'
Set-PlainCell $ws 268 6 'Active'

Set-TextCell $ws 269 1 '2026-02-06'
Set-TextCell $ws 269 2 '10:12:23'
Set-TextCell $ws 269 3 '10:00'
Set-PlainCell $ws 269 4 'Bathroom'
Set-TextCell $ws 269 5 '68.1%'
Set-PlainCell $ws 269 6 'Active'

Set-TextCell $ws 270 1 '2026-02-06'
Set-TextCell $ws 270 2 '10:12:28'
Set-TextCell $ws 270 3 '10:00'
Set-PlainCell $ws 270 4 'Bathroom'
Set-TextCell $ws 270 5 '67.4%'
Set-PlainCell $ws 270 6 'Active'

# Temperature sheet - append rows 260-270
$ws = $wb.Worksheets.Item("Temperature")
Set-TextCell $ws 260 1 '2026-02-06'
Set-TextCell $ws 260 2 '10:11:32'
Set-TextCell $ws 260 3 '10:00'
Set-PlainCell $ws 260 4 'Bathroom'
Set-TextCell $ws 260 5 '28.0C'
Set-PlainCell $ws 260 6 'Active'

Set-TextCell $ws 261 1 '2026-02-06'
Set-TextCell $ws 261 2 '10:11:35'
Set-TextCell $ws 261 3 '10:00'
Set-PlainCell $ws 261 4 'Bathroom'
Set-TextCell $ws 261 5 '28.0C'
Set-PlainCell $ws 261 6 'Active'

Set-TextCell $ws 262 1 '2026-02-06'
Set-TextCell $ws 262 2 '10:11:38'
Set-TextCell $ws 262 3 '10:00'
Set-PlainCell $ws 262 4 'Bathroom'
Set-TextCell $ws 262 5 '28.0C'
Set-PlainCell $ws 262 6 'Active'

Set-TextCell $ws 263 1 '2026-02-06'
Set-TextCell $ws 263 2 '10:11:43'
Set-TextCell $ws 263 3 '10:00'
Set-PlainCell $ws 263 4 'Bathroom'
Set-TextCell $ws 263 5 '28.0C'
Set-PlainCell $ws 263 6 'Active'

Set-TextCell $ws 264 1 '2026-02-06'
Set-TextCell $ws 264 2 '10:11:49'
Set-TextCell $ws 264 3 '10:00'
Set-PlainCell $ws 264 4 'Bathroom'
Set-TextCell $ws 264 5 '28.0C'
Set-PlainCell $ws 264 6 'Active'

Set-TextCell $ws 265 1 '2026-02-06'
Set-TextCell $ws 265 2 '10:11:54'
Set-TextCell $ws 265 3 '10:00'
Set-PlainCell $ws 265 4 'Bathroom'
Set-TextCell $ws 265 5 '27.9C'
Set-PlainCell $ws 265 6 'Active'

Set-TextCell $ws 266 1 '2026-02-06'
Set-TextCell $ws 266 2 '10:11:58'
Set-TextCell $ws 266 3 '10:00'
Set-PlainCell $ws 266 4 'Bathroom'
Set-TextCell $ws 266 5 '28.0C'
Set-PlainCell $ws 266 6 'Active'

Set-TextCell $ws 267 1 '2026-02-06'
Set-TextCell $ws 267 2 '10:12:09'
Set-TextCell $ws 267 3 '10:00'
Set-PlainCell $ws 267 4 'Bathroom'
Set-TextCell $ws 267 5 '27.9C'
Set-PlainCell $ws 267 6 'Active'

Set-TextCell $ws 268 1 '2026-02-06'
Set-TextCell $ws 268 2 '10:12:14'
Set-TextCell $ws 268 3 '10:00'
Set-PlainCell $ws 268 4 'Bathroom'
Set-TextCell $ws 268 5 '28.0C'
Set-PlainCell $ws 268 6 'Active'

Set-TextCell $ws 269 1 '2026-02-06'
Set-TextCell $ws 269 2 '10:12:24'
Set-TextCell $ws 269 3 '10:00'
Set-PlainCell $ws 269 4 'Bathroom'
Set-TextCell $ws 269 5 '28.0C'
Set-PlainCell $ws 269 6 'Active'

Set-TextCell $ws 270 1 '2026-02-06'
Set-TextCell $ws 270 2 '10:12:29'
Set-TextCell $ws 270 3 '10:00'
Set-PlainCell $ws 270 4 'Bathroom'
Set-TextCell $ws 270 5 '27.9C'
Set-PlainCell $ws 270 6 'Active'
